$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = -1.4575931196440592
$ws.Range("C2").Value = 13.093917949802517
$ws.Range("D2").Value = 18.728076225441782
$ws.Range("E2").Value = 26.223132047116678

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = -21.48109036200816
$ws.Range("C3").Value = 11.493826343276282
$ws.Range("D3").Value = 40.481075397895438
$ws.Range("E3").Value = 17.00445050643782

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
